$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "Chile"
$ws.Cells.Item(2, 2).NumberFormat = "@"
$ws.Cells.Item(2, 2).Value = "3"
$ws.Cells.Item(2, 3).Value = "Recreation"
$ws.Cells.Item(2, 4).Value = -0.07730000000000001
$ws.Cells.Item(2, 7).Value = -0.0611844452372458
$ws.Cells.Item(2, 8).Value = -0.0611844452372458
$ws.Cells.Item(2, 9).Value = -0.3740635033892258
$ws.Cells.Item(2, 10).Value = -0.3419349449852482
$ws.Cells.Item(2, 11).Value = -2.25
$ws.Cells.Item(2, 12).Value = -0.2006778451658937
$ws.Cells.Item(2, 13).Value = 0
$ws.Cells.Item(2, 14).Value = 0
$ws.Cells.Item(2, 15).Value = 0
$ws.Cells.Item(2, 16).Value = 0
$ws.Cells.Item(2, 17).Value = 0
$ws.Cells.Item(2, 18).Value = 0
$ws.Cells.Item(2, 19).Value = 0
$ws.Cells.Item(2, 21).Value = 2.92
$ws.Cells.Item(2, 22).Value = 0.03752731011438119
$ws.Cells.Item(2, 23).Value = -0.0462
$ws.Cells.Item(2, 24).Value = 0.0686033054893818
$ws.Cells.Item(2, 25).Value = -0.1148033054893818
$ws.Cells.Item(2, 26).Value = 0.1629769605349226
$ws.Cells.Item(2, 27).Value = -0.1029023746701847
$ws.Cells.Item(2, 28).Value = 0.0589379995332074
$ws.Cells.Item(2, 29).Value = -0.1604540404374995
$ws.Cells.Item(2, 30).Value = 16.47
$ws.Cells.Item(2, 31).Value = 0
$ws.Cells.Item(2, 32).Value = 16.47
$ws.Cells.Item(2, 33).Value = 13.55
$ws.Cells.Item(2, 34).Value = 0.1746924056003394
$ws.Cells.Item(2, 35).Value = 0.2676308092297692
$ws.Cells.Item(2, 36).Value = 0.1483143607705779
$ws.Cells.Item(2, 37).Value = 0.2311497782326851
$ws.Cells.Item(2, 38).Value = 0.44
$ws.Cells.Item(2, 39).Value = 0.427
$ws.Cells.Item(2, 40).Value = -4.947431661159507
$ws.Cells.Item(2, 41).Value = -9.531818181818181
$ws.Cells.Item(2, 42).Value = -4.07029137879243
$ws.Cells.Item(2, 43).Value = -9.822014051522249

# Row 3
$ws.Cells.Item(3, 1).Value = "Chile"
$ws.Cells.Item(3, 2).Value = "Club Hipico de Santiago S.A. (SNSE:HIPICO)"
$ws.Cells.Item(3, 3).Value = "Recreation"
$ws.Cells.Item(3, 4).Value = -0.07730000000000001
$ws.Cells.Item(3, 7).Value = 0.00163265306122449
$ws.Cells.Item(3, 8).Value = 0.00163265306122449
$ws.Cells.Item(3, 9).Value = -0.2744897959183673
$ws.Cells.Item(3, 10).Value = -0.2744897959183673
$ws.Cells.Item(3, 11).Value = -2.31
$ws.Cells.Item(3, 12).Value = -0.2357142857142857
$ws.Cells.Item(3, 13).Value = 0
$ws.Cells.Item(3, 14).Value = 0
$ws.Cells.Item(3, 15).Value = 0
$ws.Cells.Item(3, 16).Value = 0
$ws.Cells.Item(3, 17).Value = 0
$ws.Cells.Item(3, 18).Value = 0
$ws.Cells.Item(3, 19).Value = 0
$ws.Cells.Item(3, 21).Value = 2.45
$ws.Cells.Item(3, 22).Value = 0.08032786885245902
$ws.Cells.Item(3, 23).Value = -0.0462
$ws.Cells.Item(3, 24).Value = 0.0686033054893818
$ws.Cells.Item(3, 25).Value = -0.1148033054893818
$ws.Cells.Item(3, 26).Value = 0.1706573791902481
$ws.Cells.Item(3, 27).Value = -0.04684370918589464
$ws.Cells.Item(3, 28).Value = 0.0589379995332074
$ws.Cells.Item(3, 29).Value = -0.105781708719102
$ws.Cells.Item(3, 30).Value = 7.3
$ws.Cells.Item(3, 31).Value = 0
$ws.Cells.Item(3, 32).Value = 7.3
$ws.Cells.Item(3, 33).Value = 4.85
$ws.Cells.Item(3, 34).Value = 0.1931216931216931
$ws.Cells.Item(3, 35).Value = 0.1423001949317739
$ws.Cells.Item(3, 36).Value = 0.1371994342291372
$ws.Cells.Item(3, 37).Value = 0.09928352098259978
$ws.Cells.Item(3, 38).Value = 0.128
$ws.Cells.Item(3, 39).Value = 0.117
$ws.Cells.Item(3, 40).Value = -3.526570048309179
$ws.Cells.Item(3, 41).Value = -21.015625
$ws.Cells.Item(3, 42).Value = -2.342995169082126
$ws.Cells.Item(3, 43).Value = -22.99145299145299

# Row 4
$ws.Cells.Item(4, 1).Value = "Chile"
$ws.Cells.Item(4, 2).Value = "Unión El Golf S.A. (SNSE:UNION GOLF)"
$ws.Cells.Item(4, 3).Value = "Recreation"
$ws.Cells.Item(4, 4).Value = -0.173
$ws.Cells.Item(4, 7).Value = -0.2984962406015038
$ws.Cells.Item(4, 8).Value = -0.2984962406015038
$ws.Cells.Item(4, 9).Value = -0.8796992481203006
$ws.Cells.Item(4, 10).Value = -0.8796992481203006
$ws.Cells.Item(4, 11).Value = -1.34
$ws.Cells.Item(4, 12).Value = -1.007518796992481
$ws.Cells.Item(4, 13).Value = 0
$ws.Cells.Item(4, 14).Value = 0
$ws.Cells.Item(4, 15).Value = 0
$ws.Cells.Item(4, 16).Value = 0
$ws.Cells.Item(4, 17).Value = 0
$ws.Cells.Item(4, 18).Value = 0
$ws.Cells.Item(4, 19).Value = 0
$ws.Cells.Item(4, 21).Value = 0.47
$ws.Cells.Item(4, 22).Value = 0.05334846765039727
$ws.Cells.Item(4, 23).Value = -0.5095057034220533
$ws.Cells.Item(4, 24).Value = 0.09814125928231206
$ws.Cells.Item(4, 25).Value = -0.6076469627043654
$ws.Cells.Item(4, 26).Value = 0.1169744942832014
$ws.Cells.Item(4, 27).Value = -0.1029023746701847
$ws.Cells.Item(4, 28).Value = 0.05755166576731478
$ws.Cells.Item(4, 29).Value = -0.1604540404374995
$ws.Cells.Item(4, 30).Value = 9.17
$ws.Cells.Item(4, 31).Value = 0
$ws.Cells.Item(4, 32).Value = 9.17
$ws.Cells.Item(4, 33).Value = 8.7
$ws.Cells.Item(4, 34).Value = 0.5100111234705228
$ws.Cells.Item(4, 35).Value = 0.8955078125
$ws.Cells.Item(4, 36).Value = 0.4968589377498572
$ws.Cells.Item(4, 37).Value = 0.8904810644831115
$ws.Cells.Item(4, 38).Value = 0.219
$ws.Cells.Item(4, 39).Value = 0.219
$ws.Cells.Item(4, 40).Value = -9.61215932914046
$ws.Cells.Item(4, 41).Value = -5.342465753424658
$ws.Cells.Item(4, 42).Value = -9.119496855345911
$ws.Cells.Item(4, 43).Value = -5.342465753424658

# Row 5
$ws.Cells.Item(5, 1).Value = "Chile"
$ws.Cells.Item(5, 2).Value = "Club de Polo y Equitación San Cristóbal S.A. (SNSE:POLO)"
$ws.Cells.Item(5, 3).Value = "Recreation"
$ws.Cells.Item(5, 4).Value = 0.0916
$ws.Cells.Item(5, 7).Value = -3.719512195121951
$ws.Cells.Item(5, 8).Value = -3.719512195121951
$ws.Cells.Item(5, 9).Value = -4.073170731707317
$ws.Cells.Item(5, 10).Value = -3.023628855336172
$ws.Cells.Item(5, 11).Value = 1.4
$ws.Cells.Item(5, 12).Value = 17.07317073170731
$ws.Cells.Item(5, 13).Value = 0
$ws.Cells.Item(5, 14).Value = 0
$ws.Cells.Item(5, 15).Value = 0
$ws.Cells.Item(5, 16).Value = 0
$ws.Cells.Item(5, 17).Value = 0
$ws.Cells.Item(5, 18).Value = 0
$ws.Cells.Item(5, 19).Value = 0
$ws.Cells.Item(5, 21).Value = 0
$ws.Cells.Item(5, 22).Value = 0
$ws.Cells.Item(5, 24).Value = 0.05978287190337842
$ws.Cells.Item(5, 28).Value = 0.05978287190337842
$ws.Cells.Item(5, 30).Value = 0
$ws.Cells.Item(5, 31).Value = 0
$ws.Cells.Item(5, 32).Value = 0
$ws.Cells.Item(5, 33).Value = 0
$ws.Cells.Item(5, 34).Value = 0
$ws.Cells.Item(5, 36).Value = 0
$ws.Cells.Item(5, 38).Value = 0.093
$ws.Cells.Item(5, 39).Value = 0.091
$ws.Cells.Item(5, 40).Value = 0
$ws.Cells.Item(5, 41).Value = -3.591397849462366
$ws.Cells.Item(5, 42).Value = 0
$ws.Cells.Item(5, 43).Value = -3.670329670329671
